$d = $word.ActiveDocument

# Locate the "Account page (20" portion of the "Account page (20 hours)" run.
# We split it so that a new "_GoBack" bookmark sits right after the "20",
# matching the target structure where typing last occurred.
$findRange = $d.Content
$findRange.Find.Execute("Account page (20", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $d.Range($findRange.End, $findRange.End)

# Remove the "_GoBack" bookmark from its current (stale) location, if present.
$existing = $d.Bookmarks("_GoBack")
$existing.Delete()

# Re-add "_GoBack" right after "Account page (20" - this is what naturally
# splits the run into "Account page (20" + " hours)" once the bookmark is
# serialized between them.
$d.Bookmarks.Add("_GoBack", $insertPoint)
